$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Basic Metadata Template")

# --- Add six new resource rows (4-9) to the Basic Metadata Template sheet ---
$ws.Range("B4").Value = 'Geologynet'
$ws.Range("C4").Value = 'Collection'
$ws.Range("D4").Value = 'general geology'
$ws.Range("E4").Value = 'Website for mineral database and petrographic database'
$ws.Range("F4").Value = 'https://geologynet.com/dbases.htm'

$ws.Range("B5").Value = 'USGS National Geologic Map Database'
$ws.Range("C5").Value = 'Organization'
$ws.Range("D5").Value = 'geology maps'
$ws.Range("E5").Value = 'Website for the National Geologic map database  for the U.S.'
$ws.Range("F5").Value = 'https://ngmdb.usgs.gov/ngmdb/ngmdb_home.html'

$ws.Range("B6").Value = 'Servicio Geologico Mexicano'
$ws.Range("C6").Value = 'Organization'
$ws.Range("D6").Value = 'Geology/Geochemical maps'
$ws.Range("E6").Value = 'Website for the Mexican Geological Survy, which provides geologic, geochemical, and geophysical maps of Mexico.'
$ws.Range("F6").Value = 'http://www.gob.mx/sgm'

$ws.Range("B7").Value = 'Utah Geologic Survery'
$ws.Range("C7").Value = 'Organization'
$ws.Range("D7").Value = 'geologic data'
$ws.Range("E7").Value = 'Website provide geologic data (geology maps, sample locations, mineral resources, etc.) for the state of Utah.'
$ws.Range("F7").Value = 'http://geology.utah.gov/resources/data-databases/'

$ws.Range("B8").Value = 'RRUFF'
$ws.Range("C8").Value = 'Project'
$ws.Range("D8").Value = 'Minerals'
$ws.Range("E8").Value = 'Website providing an intensive mineral database, which has RAMAN spectra, x-ray diffraction, and chemistry data for minerals.'
$ws.Range("F8").Value = 'http://rruff.info'

$ws.Range("B9").Value = 'USGS National Map Hydrography'
$ws.Range("C9").Value = 'Organization'
$ws.Range("D9").Value = 'hydrologic '
$ws.Range("E9").Value = 'Website that provides maps about hydrologic data in the U.S. '
$ws.Range("F9").Value = 'https://viewer.nationalmap.gov/viewer/nhd.html?p=nhd'

# --- Shade the second header/example row (B2:F2) with a light grey fill ---
$ws.Range("B2:F2").Interior.Color = 14277081

# --- Update view/selection state to match the saved workbook ---
$rt = $wb.Worksheets.Item("ResourceTypes")
$rt.Activate()
$rt.Range("C40").Select()

$fl = $wb.Worksheets.Item("FieldList")
$fl.Activate()

$ws.Activate()
$ws.Range("B14").Select()
